$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# This document is a mail-merge style letter template with three placeholder
# paragraphs:
#   Para 2: "$name 2 !"       -> "Dupont Jean 2 !"
#   Para 3: "$firstname 3 ?"  -> "Dupont 3 ?"
#   Para 4: "$lastname 4 ."   -> "Jean 4 ."
#
# Each paragraph is rewritten by filling the merge value into the first run
# and deleting the now-redundant placeholder-word / trailing-literal runs
# (precise range deletes, so bookmarks / proofErr markers sitting between the
# runs are left in place rather than being swallowed by the edit).
#
# Paragraphs are processed from the end of the document back to the start so
# earlier paragraphs' character offsets are never disturbed by a later edit.
# ---------------------------------------------------------------------------

# --- Paragraph 4: "$lastname" + " 4 ." -> single run "Jean 4 ." -------------
# No bookmark / proofErr sits between these two runs, so a simple Find that
# spans both merges them cleanly into one run with the final text.
$d.Content.Find.Execute("`$lastname 4", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Jean 4", 2)

# --- Paragraph 3: "$" + "firstname" + " 3 ?" -> "Dupont 3 ?" ----------------
$p3 = $d.Paragraphs(3).Range
$p3Start = $p3.Start
$d.Range($p3Start + 10, $p3Start + 14).Delete()   # delete " 3 ?"
$d.Range($p3Start + 1,  $p3Start + 10).Delete()   # delete "firstname"
$d.Range($p3Start,      $p3Start + 1).Text = "Dupont 3 ?"

# --- Paragraph 2: "$" + "name" + " " + [bookmark] + "2 !" -> "Dupont Jean 2 !"
$p2 = $d.Paragraphs(2).Range
$p2Start = $p2.Start
$d.Range($p2Start + 6, $p2Start + 9).Delete()     # delete "2 !"
$d.Range($p2Start + 5, $p2Start + 6).Delete()     # delete " "
$d.Range($p2Start + 1, $p2Start + 5).Delete()     # delete "name"
$d.Range($p2Start,     $p2Start + 1).Text = "Dupont Jean 2 !"

Write-Output $d.Content.Text
